# "slides up to lesson 20"
# Updates the STAT 2430 Data-Viz course schedule:
#  - lesson 6 (row 20) moves out a week (week 5 -> week 6), which cascades
#    the shared date formula in column C automatically
#  - slide decks are now published for lessons 13-18b (rows 21,22,24,25,26,28,30,31)
#  - the "finding data" / "reproducible reports" slide links are renumbered
#    from 18a/18b to 19/20
#  - selection follows the edit to the last touched cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: lesson moves from week 5 to week 6 (date formula recalculates itself)
$ws.Cells.Item(20, 1).Value = 6

# Newly-published slide decks ("Slides" in column I) for lessons 13-18b
$ws.Cells.Item(21, 9).Value = "Slides"
$ws.Cells.Item(22, 9).Value = "Slides"
$ws.Cells.Item(24, 9).Value = "Slides"
$ws.Cells.Item(25, 9).Value = "Slides"
$ws.Cells.Item(26, 9).Value = "Slides"
$ws.Cells.Item(28, 9).Value = "Slides"
$ws.Cells.Item(30, 9).Value = "Slides"
$ws.Cells.Item(31, 9).Value = "Slides"

# Renumber the slide filenames for "Finding data" and "Reproducible reports"
# from the old 18a/18b naming scheme to 19/20
$ws.Cells.Item(30, 10).Value = "19-finding-data.html"
$ws.Cells.Item(31, 10).Value = "20-reproducible-reports.html"

# Move the active selection to reflect where the edits left off
$ws.Range("J31").Select()
